# Update "想去人数" (column F) counts across the four sheets to match the
# newly generated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value = 3198
$ws1.Range("F14").Value = 5788
$ws1.Range("F20").Value = 479
$ws1.Range("F21").Value = 1262
$ws1.Range("F24").Value = 2015

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 1172
$ws2.Range("F8").Value = 341
$ws2.Range("F32").Value = 188

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 1797
$ws3.Range("F5").Value = 2573
$ws3.Range("F6").Value = 1125
$ws3.Range("F10").Value = 407
$ws3.Range("F13").Value = 560

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1797
$ws4.Range("F5").Value = 2573
$ws4.Range("F6").Value = 1125
$ws4.Range("F8").Value = 407
$ws4.Range("F14").Value = 3198
$ws4.Range("F19").Value = 560
$ws4.Range("F21").Value = 341
$ws4.Range("F32").Value = 479
$ws4.Range("F40").Value = 1262
$ws4.Range("F44").Value = 2015
$ws4.Range("F45").Value = 188

$wb.Save()
